# Auto-generated edit script: refresh market-price-derived columns (H-N)
# across all 8 job sheets, per the scheduled market-data runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 26.88889
$ws.Range("I11").Value = 26.88889
$ws.Range("K11").Value = 26.88889
$ws.Range("M11").Value = 113.11111
$ws.Range("H43").Value = 2710.9375
$ws.Range("I43").Value = 2803.75
$ws.Range("J43").Value = 2432.5
$ws.Range("K43").Value = 2803.75
$ws.Range("L43").Value = 2432.5
$ws.Range("M43").Value = -2734.75
$ws.Range("N43").Value = -2570.5
$ws.Range("H45").Value = 300
$ws.Range("I45").Value = 300
$ws.Range("K45").Value = 900
$ws.Range("M45").Value = -708
$ws.Range("H69").Value = 9058.929
$ws.Range("J69").Value = 9702.083000000001
$ws.Range("L69").Value = 29106.249
$ws.Range("N69").Value = -30854.249
$ws.Range("H72").Value = 9058.929
$ws.Range("J72").Value = 9702.083000000001
$ws.Range("L72").Value = 87318.747
$ws.Range("N72").Value = -96054.747
$ws.Range("H98").Value = 2314.1428
$ws.Range("I98").Value = 2314.1428
$ws.Range("K98").Value = 2314.1428
$ws.Range("M98").Value = -816.1428000000001
$ws.Range("H106").Value = 2429.1177
$ws.Range("I106").Value = 2207.5715
$ws.Range("K106").Value = 2207.5715
$ws.Range("M106").Value = -1576.5715
$ws.Range("H107").Value = 462
$ws.Range("J107").Value = 546.4
$ws.Range("L107").Value = 546.4
$ws.Range("N107").Value = -4386.4
$ws.Range("H112").Value = 7479.5454
$ws.Range("J112").Value = 7479.5454
$ws.Range("L112").Value = 22438.6362
$ws.Range("N112").Value = -24654.6362
$ws.Range("H122").Value = 2314.1428
$ws.Range("I122").Value = 2314.1428
$ws.Range("K122").Value = 6942.428400000001
$ws.Range("M122").Value = -4492.428400000001
$ws.Range("H132").Value = 79538.63
$ws.Range("I132").Value = 79538.63
$ws.Range("K132").Value = 238615.89
$ws.Range("M132").Value = -236085.89
$ws.Range("H135").Value = 1533.7858
$ws.Range("I135").Value = 522.381
$ws.Range("J135").Value = 4568
$ws.Range("K135").Value = 4701.429
$ws.Range("L135").Value = 41112
$ws.Range("M135").Value = -2166.429
$ws.Range("N135").Value = -46182
$ws.Range("H137").Value = 2323667
$ws.Range("I137").Value = 5719.4
$ws.Range("J137").Value = 3611415.8
$ws.Range("K137").Value = 17158.2
$ws.Range("L137").Value = 10834247.4
$ws.Range("M137").Value = -14608.2
$ws.Range("N137").Value = -10839347.4
$ws.Range("H138").Value = 1708.0769
$ws.Range("I138").Value = 1039.6666
$ws.Range("K138").Value = 3118.9998
$ws.Range("M138").Value = 2021.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6292410
$ws.Range("I32").Value = 6539155.5
$ws.Range("J32").Value = 394.5
$ws.Range("K32").Value = 6539155.5
$ws.Range("L32").Value = 394.5
$ws.Range("M32").Value = -6538868.5
$ws.Range("N32").Value = -968.5
$ws.Range("H41").Value = 1633.75
$ws.Range("I41").Value = 1178.3334
$ws.Range("J41").Value = 3000
$ws.Range("K41").Value = 1178.3334
$ws.Range("L41").Value = 3000
$ws.Range("M41").Value = -764.3334
$ws.Range("N41").Value = -3828
$ws.Range("H45").Value = 2688.4
$ws.Range("I45").Value = 2497.8572
$ws.Range("J45").Value = 3133
$ws.Range("K45").Value = 2497.8572
$ws.Range("L45").Value = 3133
$ws.Range("M45").Value = -2120.8572
$ws.Range("N45").Value = -3887
$ws.Range("H74").Value = 3440.48
$ws.Range("I74").Value = 2172
$ws.Range("K74").Value = 2172
$ws.Range("M74").Value = -1298
$ws.Range("H77").Value = 3440.48
$ws.Range("I77").Value = 2172
$ws.Range("K77").Value = 10860
$ws.Range("M77").Value = -6492
$ws.Range("H122").Value = 1837
$ws.Range("I122").Value = 1712
$ws.Range("J122").Value = 1899.5
$ws.Range("K122").Value = 5136
$ws.Range("L122").Value = 5698.5
$ws.Range("M122").Value = -2686
$ws.Range("N122").Value = -10598.5
$ws.Range("H132").Value = 1484692.8
$ws.Range("I132").Value = 1696220.2
$ws.Range("K132").Value = 5088660.6
$ws.Range("M132").Value = -5086130.6
$ws.Range("H134").Value = 87775
$ws.Range("J134").Value = 87775
$ws.Range("L134").Value = 87775
$ws.Range("N134").Value = -97915

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 967.7778
$ws.Range("I94").Value = 747.3333
$ws.Range("J94").Value = 2070
$ws.Range("K94").Value = 747.3333
$ws.Range("L94").Value = 2070
$ws.Range("M94").Value = -296.3333
$ws.Range("N94").Value = -2972
$ws.Range("H99").Value = 38784.47
$ws.Range("I99").Value = 49003.363
$ws.Range("K99").Value = 49003.363
$ws.Range("M99").Value = -47505.363

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 939.8
$ws.Range("I16").Value = 966.6667
$ws.Range("J16").Value = 899.5
$ws.Range("K16").Value = 966.6667
$ws.Range("L16").Value = 899.5
$ws.Range("M16").Value = -679.6667
$ws.Range("N16").Value = -1473.5
$ws.Range("H19").Value = 3010.5715
$ws.Range("I19").Value = 4064.8
$ws.Range("J19").Value = 375
$ws.Range("K19").Value = 4064.8
$ws.Range("L19").Value = 375
$ws.Range("M19").Value = -3894.8
$ws.Range("N19").Value = -715
$ws.Range("H24").Value = 3010.5715
$ws.Range("I24").Value = 4064.8
$ws.Range("J24").Value = 375
$ws.Range("K24").Value = 4064.8
$ws.Range("L24").Value = 375
$ws.Range("M24").Value = -3894.8
$ws.Range("N24").Value = -715
$ws.Range("H31").Value = 175587.9
$ws.Range("I31").Value = 313611.8
$ws.Range("K31").Value = 313611.8
$ws.Range("M31").Value = -313316.8
$ws.Range("H34").Value = 175587.9
$ws.Range("I34").Value = 313611.8
$ws.Range("K34").Value = 313611.8
$ws.Range("M34").Value = -313409.8
$ws.Range("H105").Value = 46181.75
$ws.Range("J105").Value = 5499
$ws.Range("L105").Value = 5499
$ws.Range("N105").Value = -8993
$ws.Range("H113").Value = 939.8
$ws.Range("I113").Value = 966.6667
$ws.Range("J113").Value = 899.5
$ws.Range("K113").Value = 966.6667
$ws.Range("L113").Value = 899.5
$ws.Range("M113").Value = 1203.3333
$ws.Range("N113").Value = -5239.5
$ws.Range("H134").Value = 6573.1377
$ws.Range("I134").Value = 6874.4443
$ws.Range("K134").Value = 20623.3329
$ws.Range("M134").Value = -18088.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 765.5
$ws.Range("I36").Value = 800
$ws.Range("J36").Value = 748.25
$ws.Range("K36").Value = 2400
$ws.Range("L36").Value = 2244.75
$ws.Range("M36").Value = -2231
$ws.Range("N36").Value = -2582.75
$ws.Range("H37").Value = 89461.08
$ws.Range("J37").Value = 89461.08
$ws.Range("L37").Value = 268383.24
$ws.Range("N37").Value = -268607.24
$ws.Range("H47").Value = 2333
$ws.Range("I47").Value = 1999.5
$ws.Range("K47").Value = 5998.5
$ws.Range("M47").Value = -5567.5
$ws.Range("H134").Value = 1598.8889
$ws.Range("I134").Value = 1598.8889
$ws.Range("K134").Value = 4796.6667
$ws.Range("M134").Value = 273.3333000000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 288014.84
$ws.Range("I80").Value = 344617
$ws.Range("K80").Value = 344617
$ws.Range("M80").Value = -343619
$ws.Range("H83").Value = 288014.84
$ws.Range("I83").Value = 344617
$ws.Range("K83").Value = 1723085
$ws.Range("M83").Value = -1718093
$ws.Range("H97").Value = 1960.3
$ws.Range("I97").Value = 1714.3125
$ws.Range("K97").Value = 1714.3125
$ws.Range("M97").Value = -1218.3125
$ws.Range("H122").Value = 8373.125
$ws.Range("I122").Value = 4387.6665
$ws.Range("J122").Value = 13497.286
$ws.Range("K122").Value = 13162.9995
$ws.Range("L122").Value = 40491.858
$ws.Range("M122").Value = -10712.9995
$ws.Range("N122").Value = -45391.858
$ws.Range("H132").Value = 56233856
$ws.Range("I132").Value = 77858216
$ws.Range("K132").Value = 233574648
$ws.Range("M132").Value = -233572118
$ws.Range("H134").Value = 33592.645
$ws.Range("J134").Value = 33592.645
$ws.Range("L134").Value = 100777.935
$ws.Range("N134").Value = -105847.935
$ws.Range("H139").Value = 200000
$ws.Range("J139").Value = 200000
$ws.Range("L139").Value = 200000
$ws.Range("N139").Value = -210280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1127.6666
$ws.Range("I16").Value = 824.0909
$ws.Range("K16").Value = 824.0909
$ws.Range("M16").Value = -654.0909
$ws.Range("H40").Value = 4615.278
$ws.Range("I40").Value = 4654.75
$ws.Range("J40").Value = 4299.5
$ws.Range("K40").Value = 4654.75
$ws.Range("L40").Value = 4299.5
$ws.Range("M40").Value = -4518.75
$ws.Range("N40").Value = -4571.5
$ws.Range("H61").Value = 2833.3333
$ws.Range("I61").Value = 2833.3333
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2833.3333
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = ""
$ws.Range("N61").Value = -2631.3333
$ws.Range("H113").Value = 2833.3333
$ws.Range("I113").Value = 2833.3333
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2833.3333
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = ""
$ws.Range("N113").Value = -663.3332999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 30285.715
$ws.Range("I136").Value = 46750
$ws.Range("J136").Value = 8333.333000000001
$ws.Range("K136").Value = 140250
$ws.Range("L136").Value = 24999.999
$ws.Range("M136").Value = -137700
$ws.Range("N136").Value = -30099.999
